$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos price/volume refresh (GitHub Actions scheduled update).
# Column D values are plain-looking numeric strings (e.g. "52.70", "35.452.86")
# that must stay stored as TEXT (matching the original inlineStr cells), so we
# write them with a leading apostrophe to force Excel's text interpretation and
# then clear the transient "quote prefix" formatting it applies, leaving the cell
# un-styled exactly like the source file.

$ws.Range('D2').Value = '''35.452.86'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +0.50%  '
$ws.Range('D3').Value = '''1.925.06'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +1.57%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = '''0.739'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +12.29%  '
$ws.Range('D6').Value = '''254.97'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +4.72%  '
$ws.Range('D8').Value = '''40.65'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -2.05%  '
$ws.Range('E9').Value = '  +4.31%  '
$ws.Range('D10').Value = '''52.70'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +4.97%  '
$ws.Range('D11').Value = '''0.0741'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +4.40%  '
$ws.Range('D12').Value = '''0.100'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +0.38%  '
$ws.Range('D13').Value = '''2.203.39'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +1.48%  '
$ws.Range('D14').Value = '''12.74'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +5.77%  '
$ws.Range('D15').Value = '''0.718'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +3.68%  '
$ws.Range('D16').Value = '''1.928.02'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +2.03%  '
$ws.Range('D17').Value = '''4.91'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +1.88%  '
$ws.Range('D18').Value = '''35.444.66'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.45%  '
$ws.Range('D19').Value = '''73.69'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +3.43%  '
$ws.Range('E20').Value = '  +2.79%  '
$ws.Range('D21').Value = '''13.04'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +4.40%  '
$ws.Range('D22').Value = '''241.93'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.37%  '
$ws.Range('D23').Value = '''5.11'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +8.39%  '
$ws.Range('E24').Value = '  -0.12%  '
$ws.Range('E25').Value = '  +2.25%  '
$ws.Range('E26').Value = '  -2.68%  '
$ws.Range('D27').Value = '''168.47'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.99%  '
$ws.Range('D28').Value = '''8.68'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +4.12%  '
$ws.Range('E29').Value = '  +7.11%  '
$ws.Range('D30').Value = '''18.93'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +4.09%  '
$ws.Range('D31').Value = '''4.131.58'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +19.54%  '
$ws.Range('D32').Value = '''4.37'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +6.40%  '
$ws.Range('D33').Value = '''1.98'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +14.53%  '
$ws.Range('E34').Value = '  +23.75%  '
$ws.Range('D35').Value = '''0.0582'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +3.79%  '
$ws.Range('D36').Value = '''4.29'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +4.76%  '
$ws.Range('E37').Value = '  -0.03%  '
$ws.Range('D38').Value = '''0.915'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -1.38%  '
$ws.Range('E39').Value = '  +0.66%  '
$ws.Range('D40').Value = '''17.32'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +9.63%  '
$ws.Range('D41').Value = '''98.61'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +10.72%  '
$ws.Range('E42').Value = '  +4.54%  '
$ws.Range('E43').Value = '  +0.86%  '
$ws.Range('D44').Value = '''0.0648'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +2.00%  '
$ws.Range('E45').Value = '  +5.44%  '
$ws.Range('D46').Value = '''1.349.19'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.88%  '
$ws.Range('E47').Value = '  +0.85%  '
$ws.Range('E48').Value = '  +0.13%  '
$ws.Range('E49').Value = '  +3.47%  '
$ws.Range('D50').Value = '''45.57'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -4.05%  '
$ws.Range('D51').Value = '''2.110.53'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +1.49%  '
